$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "63.640.43"
Set-TextValue "E2" "  +1.09%  "

# Row 3
Set-TextValue "D3" "3.280.96"
Set-TextValue "E3" "  +4.82%  "

# Row 5
Set-TextValue "D5" "603.73"
Set-TextValue "E5" "  +2.47%  "

# Row 6
Set-TextValue "D6" "141.12"
Set-TextValue "E6" "  +3.64%  "

# Row 7
Set-TextValue "E7" "  +0.10%  "

# Row 8
Set-TextValue "D8" "3.279.90"
Set-TextValue "E8" "  +4.98%  "

# Row 9
Set-TextValue "D9" "0.517"
Set-TextValue "E9" "  +0.62%  "

# Row 10
Set-TextValue "E10" "  +2.97%  "

# Row 11
Set-TextValue "D11" "5.40"
Set-TextValue "E11" "  +3.16%  "

# Row 12
Set-TextValue "D12" "0.467"
Set-TextValue "E12" "  +3.03%  "

# Row 13
Set-TextValue "D13" "0.0000245"
Set-TextValue "E13" "  +0.87%  "

# Row 14
Set-TextValue "D14" "34.37"
Set-TextValue "E14" "  +1.31%  "

# Row 15
Set-TextValue "D15" "3.826.00"
Set-TextValue "E15" "  +5.15%  "

# Row 16
Set-TextValue "E16" "  +1.20%  "

# Row 17
Set-TextValue "D17" "3.282.39"
Set-TextValue "E17" "  +5.08%  "

# Row 18
Set-TextValue "D18" "63.719.59"
Set-TextValue "E18" "  +1.15%  "

# Row 19
Set-TextValue "D19" "6.81"
Set-TextValue "E19" "  +2.87%  "

# Row 20
Set-TextValue "D20" "477.43"
Set-TextValue "E20" "  +1.53%  "

# Row 21
Set-TextValue "D21" "14.02"
Set-TextValue "E21" "  -0.25%  "

# Row 22
Set-TextValue "D22" "0.726"
Set-TextValue "E22" "  +4.40%  "

# Row 23
Set-TextValue "D23" "7.98"
Set-TextValue "E23" "  +4.47%  "

# Row 24
Set-TextValue "D24" "13.50"
Set-TextValue "E24" "  +4.76%  "

# Row 25
Set-TextValue "D25" "84.04"
Set-TextValue "E25" "  -1.40%  "

# Row 26
Set-TextValue "E26" "  +0.07%  "

# Row 27
Set-TextValue "E27" "  +2.14%  "

# Row 28
Set-TextValue "B28" "FirstDigitalUSD"
Set-TextValue "C28" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D28" "1.00"
Set-TextValue "E28" "  +0.07%  "

# Row 29
Set-TextValue "B29" "NEARProtocol"
Set-TextValue "C29" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D29" "7.27"
Set-TextValue "E29" "  +6.45%  "

# Row 30
Set-TextValue "D30" "8.08"
Set-TextValue "E30" "  +2.71%  "

# Row 31
Set-TextValue "E31" "  +3.39%  "

# Row 32
Set-TextValue "D32" "28.51"
Set-TextValue "E32" "  +7.47%  "

# Row 33
Set-TextValue "D33" "0.104"
Set-TextValue "E33" "  -3.20%  "

# Row 34
Set-TextValue "E34" "  -0.54%  "

# Row 35
Set-TextValue "E35" "  +3.27%  "

# Row 36
Set-TextValue "D36" "5.93"
Set-TextValue "E36" "  +3.53%  "

# Row 37
Set-TextValue "D37" "53.40"
Set-TextValue "E37" "  +2.80%  "

# Row 38
Set-TextValue "D38" "0.0₃0733"
Set-TextValue "E38" "  +8.06%  "

# Row 39
Set-TextValue "E39" "  +3.09%  "

# Row 40
Set-TextValue "D40" "426.73"
Set-TextValue "E40" "  +2.82%  "

# Row 41
Set-TextValue "D41" "3.048.27"
Set-TextValue "E41" "  +4.56%  "

# Row 42
Set-TextValue "D42" "8.30"
Set-TextValue "E42" "  +1.73%  "

# Row 43
Set-TextValue "D43" "2.71"
Set-TextValue "E43" "  +1.51%  "

# Row 44
Set-TextValue "E44" "  +0.33%  "

# Row 45
Set-TextValue "E45" "  +1.57%  "

# Row 46
Set-TextValue "D46" "2.17"
Set-TextValue "E46" "  +3.64%  "

# Row 47
Set-TextValue "B47" "USDe"
Set-TextValue "C47" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D47" "0.999"
Set-TextValue "E47" "  +0.01%  "

# Row 48
Set-TextValue "B48" "InjectiveProtocol"
Set-TextValue "C48" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D48" "26.07"
Set-TextValue "E48" "  +3.27%  "

# Row 49
Set-TextValue "D49" "0.114"
Set-TextValue "E49" "  +1.79%  "

# Row 50
Set-TextValue "D50" "124.63"
Set-TextValue "E50" "  +3.65%  "

# Row 51
Set-TextValue "E51" "  +1.36%  "
